$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1611.1428
$ws.Range("J19").Value = 1489
$ws.Range("L19").Value = 1489
$ws.Range("N19").Value = -1839
$ws.Range("H98").Value = 2912
$ws.Range("I98").Value = 1574.75
$ws.Range("K98").Value = 1574.75
$ws.Range("M98").Value = -76.75
$ws.Range("H112").Value = 1964770.9
$ws.Range("J112").Value = 2568200.5
$ws.Range("L112").Value = 7704601.5
$ws.Range("N112").Value = -7706817.5
$ws.Range("H122").Value = 2912
$ws.Range("I122").Value = 1574.75
$ws.Range("K122").Value = 4724.25
$ws.Range("M122").Value = -2274.25
$ws.Range("H132").Value = 2701.9614
$ws.Range("I132").Value = 2738.2
$ws.Range("J132").Value = 1796
$ws.Range("K132").Value = 8214.599999999999
$ws.Range("L132").Value = 5388
$ws.Range("M132").Value = -5684.599999999999
$ws.Range("N132").Value = -10448

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35388220
$ws.Range("I32").Value = 38584844
$ws.Range("K32").Value = 38584844
$ws.Range("M32").Value = -38584557
$ws.Range("H33").Value = 29026
$ws.Range("I33").Value = 29026
$ws.Range("K33").Value = 29026
$ws.Range("M33").Value = -28697
$ws.Range("H45").Value = 2981.6155
$ws.Range("J45").Value = 3699.889
$ws.Range("L45").Value = 3699.889
$ws.Range("N45").Value = -4453.889
$ws.Range("H63").Value = 4830
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("M63").Value = -1314
$ws.Range("H66").Value = 4830
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = -6568
$ws.Range("H132").Value = 3151.8518
$ws.Range("I132").Value = 2512.4
$ws.Range("J132").Value = 4978.857
$ws.Range("K132").Value = 7537.200000000001
$ws.Range("L132").Value = 14936.571
$ws.Range("M132").Value = -5007.200000000001
$ws.Range("N132").Value = -19996.571

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 451.8889
$ws.Range("J80").Value = 434.25
$ws.Range("L80").Value = 434.25
$ws.Range("N80").Value = -2430.25
$ws.Range("H83").Value = 451.8889
$ws.Range("J83").Value = 434.25
$ws.Range("L83").Value = 2171.25
$ws.Range("N83").Value = -12155.25
$ws.Range("H94").Value = 1366.909
$ws.Range("I94").Value = 885.25
$ws.Range("J94").Value = 1944.9
$ws.Range("K94").Value = 885.25
$ws.Range("L94").Value = 1944.9
$ws.Range("M94").Value = -434.25
$ws.Range("N94").Value = -2846.9
$ws.Range("H107").Value = 2082.6
$ws.Range("J107").Value = 2671
$ws.Range("L107").Value = 2671
$ws.Range("N107").Value = -6511
$ws.Range("H134").Value = 2235506.2
$ws.Range("I134").Value = 2859483.5
$ws.Range("K134").Value = 8578450.5
$ws.Range("M134").Value = -8575915.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5761.8423
$ws.Range("I31").Value = 2592.182
$ws.Range("K31").Value = 2592.182
$ws.Range("M31").Value = -2297.182
$ws.Range("H34").Value = 5761.8423
$ws.Range("I34").Value = 2592.182
$ws.Range("K34").Value = 2592.182
$ws.Range("M34").Value = -2390.182
$ws.Range("H58").Value = 2889.575
$ws.Range("I58").Value = 2680.1714
$ws.Range("K58").Value = 2680.1714
$ws.Range("M58").Value = -2477.1714
$ws.Range("H62").Value = 3720.0833
$ws.Range("I62").Value = 3380.375
$ws.Range("J62").Value = 4399.5
$ws.Range("K62").Value = 3380.375
$ws.Range("L62").Value = 4399.5
$ws.Range("M62").Value = -2756.375
$ws.Range("N62").Value = -5647.5
$ws.Range("H65").Value = 3720.0833
$ws.Range("I65").Value = 3380.375
$ws.Range("J65").Value = 4399.5
$ws.Range("K65").Value = 16901.875
$ws.Range("L65").Value = 21997.5
$ws.Range("M65").Value = -13781.875
$ws.Range("N65").Value = -28237.5
$ws.Range("H86").Value = 38049.4
$ws.Range("J86").Value = 40247.074
$ws.Range("L86").Value = 40247.074
$ws.Range("N86").Value = -42493.074
$ws.Range("H89").Value = 38049.4
$ws.Range("J89").Value = 40247.074
$ws.Range("L89").Value = 201235.37
$ws.Range("N89").Value = -212467.37
$ws.Range("H105").Value = 1646.4445
$ws.Range("I105").Value = 1116.8572
$ws.Range("K105").Value = 1116.8572
$ws.Range("M105").Value = 630.1428000000001
$ws.Range("H107").Value = 1170.6786
$ws.Range("I107").Value = 562.73334
$ws.Range("K107").Value = 562.73334
$ws.Range("M107").Value = 1357.26666
$ws.Range("H132").Value = 3000.3872
$ws.Range("I132").Value = 2785.8147
$ws.Range("K132").Value = 8357.444100000001
$ws.Range("M132").Value = -5827.444100000001
$ws.Range("H134").Value = 3332.1428
$ws.Range("I134").Value = 3054.1667
$ws.Range("K134").Value = 9162.500100000001
$ws.Range("M134").Value = -6627.500100000001
$ws.Range("H136").Value = 2889.575
$ws.Range("I136").Value = 2680.1714
$ws.Range("K136").Value = 8040.514200000001
$ws.Range("M136").Value = -5490.514200000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1918.4615
$ws.Range("I5").Value = 1244.1
$ws.Range("K5").Value = 3732.3
$ws.Range("M5").Value = -3620.3
$ws.Range("H18").Value = 4952.25
$ws.Range("I18").Value = 3904.75
$ws.Range("K18").Value = 11714.25
$ws.Range("M18").Value = -11545.25
$ws.Range("H23").Value = 178.88235
$ws.Range("J23").Value = 218.83333
$ws.Range("L23").Value = 656.49999
$ws.Range("N23").Value = -1126.49999
$ws.Range("H97").Value = 549
$ws.Range("I97").Value = 599
$ws.Range("J97").Value = 299
$ws.Range("K97").Value = 1797
$ws.Range("L97").Value = 897
$ws.Range("M97").Value = -1301
$ws.Range("N97").Value = -1889
$ws.Range("H113").Value = 1414.5
$ws.Range("I113").Value = 577
$ws.Range("K113").Value = 1731
$ws.Range("M113").Value = 439
$ws.Range("H127").Value = 3672.5
$ws.Range("J127").Value = 3672.5
$ws.Range("L127").Value = 11017.5
$ws.Range("N127").Value = -20937.5
$ws.Range("H130").Value = 1870
$ws.Range("I130").Value = 1870
$ws.Range("K130").Value = 5610
$ws.Range("M130").Value = -590
$ws.Range("H135").Value = 1918.4615
$ws.Range("I135").Value = 1244.1
$ws.Range("K135").Value = 11196.9
$ws.Range("M135").Value = -8661.9

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1999.6666
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("M80").Value = -2
$ws.Range("H83").Value = 1999.6666
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("M83").Value = -8
$ws.Range("H102").Value = 2865.6667
$ws.Range("I102").Value = 2865.6667
$ws.Range("K102").Value = 2865.6667
$ws.Range("M102").Value = -1243.6667
$ws.Range("H132").Value = 2519.4358
$ws.Range("I132").Value = 2223.276
$ws.Range("K132").Value = 6669.828
$ws.Range("M132").Value = -4139.828

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1586.091
$ws.Range("J22").Value = 2432.8333
$ws.Range("L22").Value = 2432.8333
$ws.Range("N22").Value = -3022.8333
$ws.Range("H27").Value = 1586.091
$ws.Range("J27").Value = 2432.8333
$ws.Range("L27").Value = 2432.8333
$ws.Range("N27").Value = -2646.8333
$ws.Range("H40").Value = 41672044
$ws.Range("I40").Value = 83336710
$ws.Range("J40").Value = 7375
$ws.Range("K40").Value = 83336710
$ws.Range("L40").Value = 7375
$ws.Range("M40").Value = -83336574
$ws.Range("N40").Value = -7647
$ws.Range("H93").Value = 142860060
$ws.Range("I93").Value = 333334660
$ws.Range("J93").Value = 4100
$ws.Range("K93").Value = 333334660
$ws.Range("L93").Value = 4100
$ws.Range("M93").Value = -333333412
$ws.Range("N93").Value = -6596
$ws.Range("H132").Value = 3320.9524
$ws.Range("I132").Value = 2809.1765
$ws.Range("K132").Value = 8427.529500000001
$ws.Range("M132").Value = -5897.529500000001
$ws.Range("H136").Value = 2320.8333
$ws.Range("I136").Value = 2114.3572
$ws.Range("J136").Value = 3043.5
$ws.Range("K136").Value = 6343.071599999999
$ws.Range("L136").Value = 9130.5
$ws.Range("M136").Value = -3793.071599999999
$ws.Range("N136").Value = -14230.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1221.5143
$ws.Range("I136").Value = 968.875
$ws.Range("K136").Value = 2906.625
$ws.Range("M136").Value = -356.625
